$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13 data
$ws.Range("A13").Value = "WGE 73"
$ws.Range("B13").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("C13").Value = "15-01-2026"
$ws.Range("D13").Value = 286962
$ws.Range("E13").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("F13").Value = 34413429360
$ws.Range("G13").Value = "DCR"
$ws.Range("H13").Value = "SBIN0003229"
$ws.Range("I13").Value = "AAAFW8862C"
$ws.Range("J13").Value = "32AAAFW8862C1Z9"
$ws.Range("K13").Value = "Nithin"
$ws.Range("L13").Value = "d24339f8-1167-4e3e-9773-65b02ed18a22"
$ws.Range("M13").Value = 32555551936
$ws.Range("N13").Value = "SBIN0001890"
$ws.Range("U13").Value = "pending"
$ws.Range("V13").Value = 300
$ws.Range("X13").Value = "Fuel for Grass cutting machine 31 DEC RPA_UNIQUE_ID : cd0c4553-8c53-494a-a979-a81dc6c1fc43"
$ws.Range("Y13").Value = "HPCL, ELATHUR"
$ws.Range("Z13").Value = "FUEL EXPENSE"
$ws.Range("AA13").Value = "managerprocurement@westernidc.com"
$ws.Range("AB13").Value = "ESTIMATION NOT MATCHED"
$ws.Range("AC13").Value = 0
$ws.Range("AD13").Value = 0
$ws.Range("AE13").Value = 0

# Row 14 data
$ws.Range("A14").Value = "WGP005"
$ws.Range("B14").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("C14").Value = "15-01-2026"
$ws.Range("D14").Value = 286962
$ws.Range("E14").Value = "Western Interior Designers & Marine Contractors"
$ws.Range("F14").Value = 34413429360
$ws.Range("G14").Value = "NEFT"
$ws.Range("H14").Value = "SBIN0003229"
$ws.Range("I14").Value = "AAAFW8862C"
$ws.Range("J14").Value = "32AAAFW8862C1Z9"
$ws.Range("K14").Value = "J.B. SALES AND SERVICE"
$ws.Range("L14").Value = "6f2bd2cd-e0f6-4007-b7a0-c83d3e7cb5b7"
$ws.Range("M14").Value = 50200049560664
$ws.Range("N14").Value = "HDFC0000072"
$ws.Range("U14").Value = "pending"
$ws.Range("V14").Value = 12897
$ws.Range("X14").Value = "Purchase of Consumables for GSL, HULL RPA_UNIQUE_ID : 1572a51e-2bc9-499b-9506-c9301e1b1d2f"
$ws.Range("Y14").Value = "GSL, HULL"
$ws.Range("Z14").Value = "SITE PURCHASE"
$ws.Range("AA14").Value = "managerprocurement@westernidc.com"
$ws.Range("AB14").Value = "ESTIMATION NOT MATCHED"
$ws.Range("AC14").Value = 0
$ws.Range("AD14").Value = 0
$ws.Range("AE14").Value = 0
